$d = $word.ActiveDocument

# 1) Mom's line when leaving for work: add "_eyes_closed" to her expression tag.
$d.Content.Find.Execute(
    "Mom (neutral smiling): I’m going to work now, so I’ll see you later. Your breakfast is on the table.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Mom (neutral smiling_eyes_closed): I’m going to work now, so I’ll see you later. Your breakfast is on the table.",
    2)

# 2) Insert a new "Mom (exit):" stage-direction paragraph right before the
#    paragraph describing her leaving the room.
$r = $d.Content
$r.Find.Execute("My mom gets up and leaves my room, and after a brief moment of hesitation I get up and follow.",
    $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$r.InsertBefore("Mom (exit):`r")

# 3) Merge the three runs that make up the umbrella/arms-touching sentence
#    into a single run.
$r2 = $d.Content
$r2.Find.Execute("My umbrella isn’t very large, so our arms touch as we walk to school. Both of us are neither brave nor embarrassed enough to point it out, and so we continue on. It’s actually kinda nice since it’s a little chilly outside, but Mara’s warm.",
    $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$r2.Text = "PLACEHOLDER_MERGE_TOKEN"
$r3 = $d.Range($r2.Start, $r2.Start + 23)
$r3.Text = "My umbrella isn’t very large, so our arms touch as we walk to school. Both of us are neither brave nor embarrassed enough to point it out, and so we continue on. It’s actually kinda nice since it’s a little chilly outside, but Mara’s warm."
